# Auto-generated edit script: updates market-board derived columns (H-N)
# across multiple Leve sheets, per the scheduled data-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 81.25
$ws.Range("I8").Value = 81.25
$ws.Range("K8").Value = 243.75
$ws.Range("M8").Value = -104.75

$ws.Range("H81").Value = 33000
$ws.Range("J81").Value = 33000
$ws.Range("L81").Value = 33000
$ws.Range("N81").Value = -34996

$ws.Range("H84").Value = 33000
$ws.Range("J84").Value = 33000
$ws.Range("L84").Value = 99000
$ws.Range("N84").Value = -108984

$ws.Range("H123").Value = 45898.332
$ws.Range("J123").Value = 48533.93
$ws.Range("L123").Value = 48533.93
$ws.Range("N123").Value = -58333.93

$ws.Range("H124").Value = 26936.363
$ws.Range("J124").Value = 26936.363
$ws.Range("L124").Value = 26936.363
$ws.Range("N124").Value = -36756.363

$ws.Range("H133").Value = 39689.25
$ws.Range("J133").Value = 39689.25
$ws.Range("L133").Value = 39689.25
$ws.Range("N133").Value = -49809.25

$ws.Range("H134").Value = 47873.855
$ws.Range("J134").Value = 47873.855
$ws.Range("L134").Value = 47873.855
$ws.Range("N134").Value = -58013.855

$ws.Range("H135").Value = 456514.97
$ws.Range("I135").Value = 456514.97
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4108634.73
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4106099.73
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 4930.769
$ws.Range("I6").Value = 4500
$ws.Range("J6").Value = 4966.6665
$ws.Range("K6").Value = 4500
$ws.Range("L6").Value = 4966.6665
$ws.Range("M6").Value = -4327
$ws.Range("N6").Value = -5312.6665

$ws.Range("H122").Value = 10578.286
$ws.Range("I122").Value = 11062.2
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 33186.60000000001
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -30736.60000000001
$ws.Range("N122").Value = -7600

$ws.Range("H123").Value = 56460
$ws.Range("J123").Value = 56460
$ws.Range("L123").Value = 56460
$ws.Range("N123").Value = -66260

$ws.Range("H133").Value = 28996
$ws.Range("J133").Value = 28996
$ws.Range("L133").Value = 28996
$ws.Range("N133").Value = -34056

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 32600
$ws.Range("J126").Value = 32600
$ws.Range("L126").Value = 32600
$ws.Range("N126").Value = -42480

$ws.Range("H130").Value = 48296
$ws.Range("J130").Value = 48296
$ws.Range("L130").Value = 48296
$ws.Range("N130").Value = -58336

$ws.Range("H132").Value = 97485.71000000001
$ws.Range("J132").Value = 97485.71000000001
$ws.Range("L132").Value = 97485.71000000001
$ws.Range("N132").Value = -107605.71

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 9800
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 9800
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 9800
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -10078

$ws.Range("H20").Value = 48949.5
$ws.Range("J20").Value = 48949.5
$ws.Range("L20").Value = 48949.5
$ws.Range("N20").Value = -49421.5

$ws.Range("H30").Value = 48949.5
$ws.Range("J30").Value = 48949.5
$ws.Range("L30").Value = 48949.5
$ws.Range("N30").Value = -49131.5

$ws.Range("H39").Value = 111134904
$ws.Range("I39").Value = 26263
$ws.Range("J39").Value = 200021810
$ws.Range("K39").Value = 26263
$ws.Range("L39").Value = 200021810
$ws.Range("M39").Value = -25872
$ws.Range("N39").Value = -200022592

$ws.Range("H49").Value = 111134904
$ws.Range("I49").Value = 26263
$ws.Range("J49").Value = 200021810
$ws.Range("K49").Value = 26263
$ws.Range("L49").Value = 200021810
$ws.Range("M49").Value = -26081
$ws.Range("N49").Value = -200022174

$ws.Range("H58").Value = 1509.037
$ws.Range("I58").Value = 1579.1177
$ws.Range("J58").Value = 1389.9
$ws.Range("K58").Value = 1579.1177
$ws.Range("L58").Value = 1389.9
$ws.Range("M58").Value = -1376.1177
$ws.Range("N58").Value = -1795.9

$ws.Range("H60").Value = 24815.143
$ws.Range("J60").Value = 24815.143
$ws.Range("L60").Value = 24815.143
$ws.Range("N60").Value = -25837.143

$ws.Range("H68").Value = 15360.4
$ws.Range("J68").Value = 15360.4
$ws.Range("L68").Value = 15360.4
$ws.Range("N68").Value = -16858.4

$ws.Range("H71").Value = 15360.4
$ws.Range("J71").Value = 15360.4
$ws.Range("L71").Value = 46081.2
$ws.Range("N71").Value = -53569.2

$ws.Range("H80").Value = 28000
$ws.Range("J80").Value = 28000
$ws.Range("L80").Value = 28000
$ws.Range("N80").Value = -30246

$ws.Range("H83").Value = 28000
$ws.Range("J83").Value = 28000
$ws.Range("L83").Value = 84000
$ws.Range("N83").Value = -95232

$ws.Range("H109").Value = 11425
$ws.Range("J109").Value = 11425
$ws.Range("L109").Value = 11425
$ws.Range("N109").Value = -13505

$ws.Range("H128").Value = 48949.5
$ws.Range("J128").Value = 48949.5
$ws.Range("L128").Value = 48949.5
$ws.Range("N128").Value = -58909.5

$ws.Range("H136").Value = 1509.037
$ws.Range("I136").Value = 1579.1177
$ws.Range("J136").Value = 1389.9
$ws.Range("K136").Value = 4737.3531
$ws.Range("L136").Value = 4169.700000000001
$ws.Range("M136").Value = -2187.3531
$ws.Range("N136").Value = -9269.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 913.88
$ws.Range("J131").Value = 933.625
$ws.Range("L131").Value = 2800.875
$ws.Range("N131").Value = -12880.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7800
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 8444.444
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 8444.444
$ws.Range("M5").Value = -1888
$ws.Range("N5").Value = -8668.444

$ws.Range("H46").Value = 9398
$ws.Range("I46").Value = 4250
$ws.Range("J46").Value = 19694
$ws.Range("K46").Value = 4250
$ws.Range("L46").Value = 19694
$ws.Range("M46").Value = -4094
$ws.Range("N46").Value = -20006

$ws.Range("H123").Value = 38453.89
$ws.Range("J123").Value = 38453.89
$ws.Range("L123").Value = 38453.89
$ws.Range("N123").Value = -43353.89

$ws.Range("H124").Value = 49092
$ws.Range("J124").Value = 49092
$ws.Range("L124").Value = 49092
$ws.Range("N124").Value = -58912

$ws.Range("H130").Value = 33043.08
$ws.Range("J130").Value = 33043.08
$ws.Range("L130").Value = 33043.08
$ws.Range("N130").Value = -43083.08

$ws.Range("H133").Value = 38588
$ws.Range("J133").Value = 38588
$ws.Range("L133").Value = 38588
$ws.Range("N133").Value = -48708

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 50820
$ws.Range("J123").Value = 50820
$ws.Range("L123").Value = 50820
$ws.Range("N123").Value = -60620

$ws.Range("H127").Value = 50749
$ws.Range("J127").Value = 50749
$ws.Range("L127").Value = 50749
$ws.Range("N127").Value = -60669

$ws.Range("H128").Value = 43326.668
$ws.Range("J128").Value = 43326.668
$ws.Range("L128").Value = 43326.668
$ws.Range("N128").Value = -53286.668

$ws.Range("H130").Value = 46740
$ws.Range("J130").Value = 46740
$ws.Range("L130").Value = 46740
$ws.Range("N130").Value = -56780

$ws.Range("H132").Value = 5775.227
$ws.Range("I132").Value = 6565.2856
$ws.Range("J132").Value = 4392.625
$ws.Range("K132").Value = 19695.8568
$ws.Range("L132").Value = 13177.875
$ws.Range("M132").Value = -17165.8568
$ws.Range("N132").Value = -18237.875

$ws.Range("H136").Value = 1870.258
$ws.Range("I136").Value = 1765.9333
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 5297.7999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2747.7999
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 24000
$ws.Range("J109").Value = 24000
$ws.Range("L109").Value = 24000
$ws.Range("N109").Value = -26774

$ws.Range("H123").Value = 27307.3
$ws.Range("J123").Value = 28296.143
$ws.Range("L123").Value = 28296.143
$ws.Range("N123").Value = -38096.143

$ws.Range("H125").Value = 47153.332
$ws.Range("J125").Value = 47153.332
$ws.Range("L125").Value = 47153.332
$ws.Range("N125").Value = -56993.332

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H135").Value = 62472.273
$ws.Range("J135").Value = 65719.5
$ws.Range("L135").Value = 65719.5
$ws.Range("N135").Value = -75859.5

